$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D2 (joining date) typo: 45211 -> 45212
$ws.Range("D2").Value = 45212

# Row 8 - new employee record (id 7)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Pranay"
$ws.Range("C8").Value = "inactive"
$ws.Range("D8").Value = 45211
$ws.Range("E8").Value = 41197
$ws.Range("F8").Value = "React"
$ws.Range("G8").Value = 20005
$ws.Range("H8").Value = "Kota"

# Row 9 - new employee record (id 8)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Pranay"
$ws.Range("C9").Value = "inactive"
$ws.Range("D9").Value = 45211
$ws.Range("E9").Value = 41197
$ws.Range("F9").Value = "React"
$ws.Range("G9").Value = 20006
$ws.Range("H9").Value = "Kota"

# Update view selection to B14 (and drop the old topLeftCell scroll position)
$ws.Range("B14").Select()
